{"js": "const replacements = [\n  [\"50\u00f75=\", \"79\u00f75=\"],\n  [\"65\u00f74=\", \"95\u00f73=\"],\n  [\"69\u00f77=\", \"80\u00f78=\"],\n  [\"88\u00f79=\", \"31\u00f75=\"],\n  [\"14\u00f73=\", \"22\u00f74=\"],\n  [\"58\u00f77=\", \"32\u00f73=\"],\n  [\"67\u00f77=\", \"23\u00f72=\"],\n  [\"21\u00f79=\", \"43\u00f73=\"],\n  [\"56\u00f76=\", \"84\u00f74=\"],\n  [\"28\u00f76=\", \"29\u00f74=\"],\n  [\"27\u00f76=\", \"60\u00f74=\"],\n  [\"60\u00f76=\", \"51\u00f78=\"],\n  [\"92\u00f79=\", \"93\u00f77=\"],\n  [\"29\u00f75=\", \"45\u00f77=\"],\n  [\"44\u00f76=\", \"99\u00f78=\"],\n  [\"25\u00f76=\", \"43\u00f75=\"],\n  [\"64\u00f79=\", \"30\u00f77=\"],\n  [\"50\u00f74=\", \"86\u00f72=\"],\n  [\"65\u00f78=\", \"59\u00f79=\"],\n  [\"54\u00f76=\", \"44\u00f77=\"],\n  [\"37\u00f74=\", \"87\u00f76=\"],\n  [\"24\u00f73=\", \"81\u00f73=\"],\n  [\"15\u00f74=\", \"92\u00f72=\"],\n  [\"35\u00f77=\", \"70\u00f72=\"],\n  [\"42\u00f75=\", \"16\u00f75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old=\"50\u00f75=\"; New=\"79\u00f75=\"},\n    @{Old=\"65\u00f74=\"; New=\"95\u00f73=\"},\n    @{Old=\"69\u00f77=\"; New=\"80\u00f78=\"},\n    @{Old=\"88\u00f79=\"; New=\"31\u00f75=\"},\n    @{Old=\"14\u00f73=\"; New=\"22\u00f74=\"},\n    @{Old=\"58\u00f77=\"; New=\"32\u00f73=\"},\n    @{Old=\"67\u00f77=\"; New=\"23\u00f72=\"},\n    @{Old=\"21\u00f79=\"; New=\"43\u00f73=\"},\n    @{Old=\"56\u00f76=\"; New=\"84\u00f74=\"},\n    @{Old=\"28\u00f76=\"; New=\"29\u00f74=\"},\n    @{Old=\"27\u00f76=\"; New=\"60\u00f74=\"},\n    @{Old=\"60\u00f76=\"; New=\"51\u00f78=\"},\n    @{Old=\"92\u00f79=\"; New=\"93\u00f77=\"},\n    @{Old=\"29\u00f75=\"; New=\"45\u00f77=\"},\n    @{Old=\"44\u00f76=\"; New=\"99\u00f78=\"},\n    @{Old=\"25\u00f76=\"; New=\"43\u00f75=\"},\n    @{Old=\"64\u00f79=\"; New=\"30\u00f77=\"},\n    @{Old=\"50\u00f74=\"; New=\"86\u00f72=\"},\n    @{Old=\"65\u00f78=\"; New=\"59\u00f79=\"},\n    @{Old=\"54\u00f76=\"; New=\"44\u00f77=\"},\n    @{Old=\"37\u00f74=\"; New=\"87\u00f76=\"},\n    @{Old=\"24\u00f73=\"; New=\"81\u00f73=\"},\n    @{Old=\"15\u00f74=\"; New=\"92\u00f72=\"},\n    @{Old=\"35\u00f77=\"; New=\"70\u00f72=\"},\n    @{Old=\"42\u00f75=\"; New=\"16\u00f75=\"}\n)\n\n$wdReplaceAll = 2\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, $wdReplaceAll)\n}\n"}
